$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before current row 187 (shifts old rows 187-233 down to 189-235)
$ws.Rows("187:188").Insert()

# New row 187: same data as the original row 187, but with an updated date (44694)
$ws.Range("A187").Value = 1
$ws.Range("B187").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C187").Value = "Arica y Parinacota"
$ws.Range("D187").Value = 44694
$ws.Range("E187").Value = 15
$ws.Range("F187").Value = "Fruta"
$ws.Range("G187").Value = 100102
$ws.Range("H187").Value = "Cítricos"
$ws.Range("I187").Value = 100102003
$ws.Range("J187").Value = "Limón"
$ws.Range("K187").Value = "Sutil De Gase"
$ws.Range("L187").Value = "Primera"
$ws.Range("M187").Value = 250
$ws.Range("N187").Value = 29000
$ws.Range("O187").Value = 30000
$ws.Range("P187").Value = 29500
$ws.Range("Q187").Value = "$/caja 24 kilos"
$ws.Range("R187").Value = "Perú"
$ws.Range("S187").Value = 1229
$ws.Range("T187").Value = 24

# New row 188: same data pattern as the original row 188, but with an updated date (44694)
# and updated price figures (N, O, P, S)
$ws.Range("A188").Value = 1
$ws.Range("B188").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C188").Value = "Arica y Parinacota"
$ws.Range("D188").Value = 44694
$ws.Range("E188").Value = 15
$ws.Range("F188").Value = "Fruta"
$ws.Range("G188").Value = 100102
$ws.Range("H188").Value = "Cítricos"
$ws.Range("I188").Value = 100102003
$ws.Range("J188").Value = "Limón"
$ws.Range("K188").Value = "Tahití"
$ws.Range("L188").Value = "Primera"
$ws.Range("M188").Value = 300
$ws.Range("N188").Value = 25000
$ws.Range("O188").Value = 26000
$ws.Range("P188").Value = 25500
$ws.Range("Q188").Value = "$/caja 24 kilos"
$ws.Range("R188").Value = "Perú"
$ws.Range("S188").Value = 1062
$ws.Range("T188").Value = 24
